$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows 4 and 5 entirely (they are removed in the new layout)
$ws.Range("A4:E5").ClearContents()

# Copy the header style (s="1", bold/centered/bordered) onto the new F1 header
# cell and onto the two numeric id cells A2:A3 before writing their values.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)

# Clear old A1 header text ("code" moves to C1) but keep the cell/style in place
$ws.Range("A1").ClearContents()

# New header row (row 1): B1..F1
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

# Row 2 data
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eng"
$ws.Range("C2").Value = "DKS"
$ws.Range("D2").Value = "Desktop"
$ws.Range("E2").Value = "Desktop Computer"
$ws.Range("F2").Value = $true

# Row 3 data
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "fra"
$ws.Range("C3").Value = "DKS"
$ws.Range("D3").Value = "Ordinateur"
$ws.Range("E3").Value = "Ordinateurs de bureau"
$ws.Range("F3").Value = $true
